$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 600
$ws.Range("J31").Value = 1000
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3460

$ws.Range("H100").Value = 11112504
$ws.Range("I100").Value = 14493537
$ws.Range("J100").Value = 3394.5715
$ws.Range("K100").Value = 14493537
$ws.Range("L100").Value = 3394.5715
$ws.Range("M100").Value = -14492996
$ws.Range("N100").Value = -4476.5715

$ws.Range("H137").Value = 1119
$ws.Range("J137").Value = 1139.2667
$ws.Range("L137").Value = 3417.800099999999
$ws.Range("N137").Value = -8517.8001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 763.75
$ws.Range("I74").Value = 707.8333
$ws.Range("J74").Value = 805.6875
$ws.Range("K74").Value = 707.8333
$ws.Range("L74").Value = 805.6875
$ws.Range("M74").Value = 166.1667
$ws.Range("N74").Value = -2553.6875

$ws.Range("H77").Value = 763.75
$ws.Range("I77").Value = 707.8333
$ws.Range("J77").Value = 805.6875
$ws.Range("K77").Value = 3539.1665
$ws.Range("L77").Value = 4028.4375
$ws.Range("M77").Value = 828.8334999999997
$ws.Range("N77").Value = -12764.4375

$ws.Range("H88").Value = 5511.1113
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 5950
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 5950
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -6762

$ws.Range("H91").Value = 5511.1113
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 5950
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 5950
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -8758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 69117.13
$ws.Range("I86").Value = 2512.5
$ws.Range("J86").Value = 335535.66
$ws.Range("K86").Value = 2512.5
$ws.Range("L86").Value = 335535.66
$ws.Range("M86").Value = -1389.5
$ws.Range("N86").Value = -337781.66

$ws.Range("H89").Value = 69117.13
$ws.Range("I89").Value = 2512.5
$ws.Range("J89").Value = 335535.66
$ws.Range("K89").Value = 12562.5
$ws.Range("L89").Value = 1677678.3
$ws.Range("M89").Value = -6946.5
$ws.Range("N89").Value = -1688910.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1023350.2
$ws.Range("I31").Value = 2175502
$ws.Range("K31").Value = 2175502
$ws.Range("M31").Value = -2175207

$ws.Range("H34").Value = 1023350.2
$ws.Range("I34").Value = 2175502
$ws.Range("K34").Value = 2175502
$ws.Range("M34").Value = -2175300

$ws.Range("H134").Value = 15522615
$ws.Range("I134").Value = 18000494
$ws.Range("J134").Value = 35878.5
$ws.Range("K134").Value = 54001482
$ws.Range("L134").Value = 107635.5
$ws.Range("M134").Value = -53998947
$ws.Range("N134").Value = -112705.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1053.5
$ws.Range("I117").Value = 809.6667
$ws.Range("J117").Value = 1419.25
$ws.Range("K117").Value = 2429.0001
$ws.Range("L117").Value = 4257.75
$ws.Range("M117").Value = 1012.9999
$ws.Range("N117").Value = -11141.75

$ws.Range("H129").Value = 16667713
$ws.Range("I129").Value = 1365
$ws.Range("J129").Value = 20834300
$ws.Range("K129").Value = 4095
$ws.Range("L129").Value = 62502900
$ws.Range("M129").Value = 905
$ws.Range("N129").Value = -62512900

$ws.Range("H131").Value = 6024993
$ws.Range("I131").Value = 611.86664
$ws.Range("J131").Value = 7353901
$ws.Range("K131").Value = 1835.59992
$ws.Range("L131").Value = 22061703
$ws.Range("M131").Value = 3204.40008
$ws.Range("N131").Value = -22071783

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -23744

$ws.Range("H95").Value = 13614.333
$ws.Range("J95").Value = 13614.333
$ws.Range("L95").Value = 13614.333
$ws.Range("N95").Value = -19106.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1362.5883
$ws.Range("I16").Value = 1091.1818
$ws.Range("J16").Value = 1860.1666
$ws.Range("K16").Value = 1091.1818
$ws.Range("L16").Value = 1860.1666
$ws.Range("M16").Value = -921.1818000000001
$ws.Range("N16").Value = -2200.1666

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H100").Value = 2269.3333
$ws.Range("I100").Value = 1704.8889
$ws.Range("J100").Value = 2833.7778
$ws.Range("K100").Value = 1704.8889
$ws.Range("L100").Value = 2833.7778
$ws.Range("M100").Value = -1163.8889
$ws.Range("N100").Value = -3915.7778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 2500
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 2500
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 2500
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -2980

$ws.Range("H22").Value = 2990
$ws.Range("J22").Value = 2990
$ws.Range("L22").Value = 2990
$ws.Range("N22").Value = -3576

$ws.Range("H24").Value = 1002003.3
$ws.Range("J24").Value = 1002003.3
$ws.Range("L24").Value = 1002003.3
$ws.Range("N24").Value = -1002463.3

$ws.Range("H25").Value = 1668942.4
$ws.Range("J25").Value = 1668942.4
$ws.Range("L25").Value = 1668942.4
$ws.Range("N25").Value = -1669528.4

$ws.Range("H28").Value = 4000
$ws.Range("J28").Value = 4000
$ws.Range("L28").Value = 4000
$ws.Range("N28").Value = -4696

$ws.Range("H31").Value = 4000
$ws.Range("J31").Value = 4000
$ws.Range("L31").Value = 4000
$ws.Range("N31").Value = -4696

$ws.Range("H100").Value = 1089.125
$ws.Range("I100").Value = 800
$ws.Range("J100").Value = 1130.4286
$ws.Range("K100").Value = 1600
$ws.Range("L100").Value = 2260.8572
$ws.Range("M100").Value = -1059
$ws.Range("N100").Value = -3342.8572

$ws.Range("H132").Value = 57695584
$ws.Range("I132").Value = 90002824
$ws.Range("J132").Value = 4085.8572
$ws.Range("K132").Value = 270008472
$ws.Range("L132").Value = 12257.5716
$ws.Range("M132").Value = -270005942
$ws.Range("N132").Value = -17317.5716
